$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.653.01'
$ws.Range("E2").Value = '  +3.93%  '

# Row 3
$ws.Range("D3").Value = '1.800.12'
$ws.Range("E3").Value = '  +0.53%  '

# Row 4
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = "'313.56"
$ws.Range("E5").Value = '  +0.10%  '

# Row 6
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").Value = "'0.5356"
$ws.Range("E7").Value = '  -0.84%  '

# Row 8
$ws.Range("D8").Value = "'0.3791"
$ws.Range("E8").Value = '  +0.62%  '

# Row 9
$ws.Range("D9").Value = "'0.07550"
$ws.Range("E9").Value = '  +0.49%  '

# Row 10
$ws.Range("D10").Value = "'42.69"
$ws.Range("E10").Value = '  -0.23%  '

# Row 11
$ws.Range("D11").Value = "'1.121"
$ws.Range("E11").Value = '  +0.78%  '

# Row 12
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = "'21.10"
$ws.Range("E12").Value = '  +1.09%  '

# Row 13
$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").Value = "'0.9992"
$ws.Range("E13").Value = '  -0.09%  '

# Row 14
$ws.Range("D14").Value = "'6.195"
$ws.Range("E14").Value = '  +0.29%  '

# Row 15
$ws.Range("D15").Value = "'7.490"

# Row 16
$ws.Range("D16").Value = '1.796.44'

# Row 17
$ws.Range("D17").Value = "'90.58"
$ws.Range("E17").Value = '  -0.09%  '

# Row 18
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = '  -0.14%  '

# Row 19
$ws.Range("D19").Value = "'0.06470"
$ws.Range("E19").Value = '  -0.31%  '

# Row 20
$ws.Range("D20").Value = "'0.9997"
$ws.Range("E20").Value = '  +0.01%  '

# Row 21
$ws.Range("D21").Value = "'17.29"
$ws.Range("E21").Value = '  +1.97%  '

# Row 22
$ws.Range("D22").Value = "'5.937"
$ws.Range("E22").Value = '  -0.01%  '

# Row 23
$ws.Range("D23").Value = '28.660.24'
$ws.Range("E23").Value = '  +3.82%  '

# Row 24
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = '  +0.10%  '

# Row 25
$ws.Range("D25").Value = "'2.097"
$ws.Range("E25").Value = '  +0.45%  '

# Row 26
$ws.Range("D26").Value = "'161.02"
$ws.Range("E26").Value = '  +3.76%  '

# Row 27
$ws.Range("D27").Value = "'20.54"
$ws.Range("E27").Value = '  +0.44%  '

# Row 28
$ws.Range("D28").Value = "'2.386"
$ws.Range("E28").Value = '  +0.41%  '

# Row 29
$ws.Range("D29").Value = '2.005.89'
$ws.Range("E29").Value = '  +0.54%  '

# Row 30
$ws.Range("D30").Value = "'123.69"
$ws.Range("E30").Value = '  +1.50%  '

# Row 31
$ws.Range("D31").Value = "'1.116"
$ws.Range("E31").Value = '  -1.03%  '

# Row 32
$ws.Range("D32").Value = "'0.1027"
$ws.Range("E32").Value = '  -0.69%  '

# Row 33
$ws.Range("D33").Value = "'5.703"
$ws.Range("E33").Value = '  +0.61%  '

# Row 34
$ws.Range("D34").Value = "'3.691"
$ws.Range("E34").Value = '  +2.07%  '

# Row 35
$ws.Range("D35").Value = "'0.2269"
$ws.Range("E35").Value = '  +8.43%  '

# Row 36
$ws.Range("D36").Value = "'0.06498"
$ws.Range("E36").Value = '  +7.85%  '

# Row 37
$ws.Range("D37").Value = "'8.942"
$ws.Range("E37").Value = '  +2.95%  '

# Row 38
$ws.Range("D38").Value = "'0.02308"
$ws.Range("E38").Value = '  +1.15%  '

# Row 39
$ws.Range("D39").Value = "'5.074"
$ws.Range("E39").Value = '  +1.61%  '

# Row 40
$ws.Range("D40").Value = "'11.46"
$ws.Range("E40").Value = '  +0.66%  '

# Row 41
$ws.Range("D41").Value = "'0.6291"
$ws.Range("E41").Value = '  +0.88%  '

# Row 42
$ws.Range("D42").Value = "'1.210"
$ws.Range("E42").Value = '  +5.67%  '

# Row 43
$ws.Range("D43").Value = "'0.9992"
$ws.Range("E43").Value = '  -0.08%  '

# Row 44
$ws.Range("D44").Value = "'1.389"
$ws.Range("E44").Value = '  -1.66%  '

# Row 45
$ws.Range("D45").Value = "'13.58"
$ws.Range("E45").Value = '  +2.30%  '

# Row 46
$ws.Range("D46").Value = "'0.5922"
$ws.Range("E46").Value = '  +1.00%  '

# Row 47
$ws.Range("D47").Value = "'3.665"
$ws.Range("E47").Value = '  +1.11%  '

# Row 48
$ws.Range("D48").Value = "'126.09"
$ws.Range("E48").Value = '  +3.60%  '

# Row 49
$ws.Range("D49").Value = "'1.976"
$ws.Range("E49").Value = '  +3.40%  '

# Row 50
$ws.Range("D50").Value = "'1.163"
$ws.Range("E50").Value = '  +2.59%  '
